$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, E, F, G, H, I, M, N across rows 2-25
# (row index => column letter => new value)
$data = @{
    2  = @{ B='11.79008979645891';  C='9.146297098557223'; E='11.64193756024586'; F='16.86991607391245'; G='25.80280529302172'; H='12.96756632760109'; I='17.99397222990934'; M='14.40239269569655'; N='16.40324474656374' }
    3  = @{ B='11.20128780790012';  C='8.617083079331648'; E='11.53414368500319'; F='15.89584955866815'; G='25.47253618121673'; H='12.99458165896553'; I='18.06806836529239'; M='14.1179314706999';  N='16.47010806177864' }
    4  = @{ B='10.82527740164572';  C='8.273478209271925'; E='11.47156272440864'; F='15.26997757108489'; G='25.28127003946212'; H='13.01474348701838'; I='18.11916865240467'; M='13.94422934883281'; N='16.51308396530331' }
    5  = @{ B='10.66860470937609';  C='8.128799319758761'; E='11.44699204141228'; F='15.008197319934';   G='25.20633087119935'; H='13.02385443183121'; I='18.14139313045144'; M='13.87378864532545'; N='16.53108177267938' }
    6  = @{ B='10.64238761962339';  C='8.104495491987191'; E='11.44296905199232'; F='14.96433081551589'; G='25.1940716026063';  H='13.02542122196905'; I='18.14516782856428'; M='13.86211585819523'; N='16.53409962101097' }
    7  = @{ B='10.8231781256704';   C='8.271545819972971'; E='11.47122755222083'; F='15.26647399323133'; G='25.28024708951274'; H='13.014862742738';   I='18.11946271983773'; M='13.94327782935429'; N='16.51332472510162' }
    8  = @{ B='11.59018521961684';  C='8.967716882844689'; E='11.60403774370236'; F='16.5399640634477';  G='25.68660404705504'; H='12.97613707579701'; I='18.01835207780978'; M='14.30417319109519'; N='16.42590135536888' }
    9  = @{ B='12.97224001865727';  C='10.18367893268525'; E='11.891876952679';   F='19.00274580682531'; G='26.56977634533969'; H='12.92872101361982'; I='17.86492929360339'; M='15.01503397610182'; N='16.26963806331387' }
    10 = @{ B='13.90551102086756';  C='10.98509641590152'; E='12.1183151996688';  F='20.67494806633232'; G='27.26394515072758'; H='12.91147990922998'; I='17.78006459007004'; M='15.53343875166049'; N='16.16397857394625' }
    11 = @{ B='14.31105181644751';  C='11.32962487861414'; E='12.22417670218902'; F='21.3917225636224';  G='27.58792955498873'; H='12.90749426047885'; I='17.74761052665655'; M='15.76728264823201'; N='16.11787565904532' }
    12 = @{ B='14.46180941325502';  C='11.45720412852704'; E='12.26463902908495'; F='21.65686569030329'; G='27.71165310671114'; H='12.90654210756537'; I='17.73621376202193'; M='15.85545835720438'; N='16.10069821199387' }
    13 = @{ B='14.42946722396115';  C='11.42985605482051'; E='12.25590864155965'; F='21.60004134736742'; G='27.68496313518158'; H='12.90672235781333'; I='17.7386284184533';  M='15.83648625431124'; N='16.1043852196585'  }
    14 = @{ B='14.32351148742788';  C='11.34017875777904'; E='12.22749825448128'; F='21.4136618050453';  G='27.59808838259336'; H='12.90740474810743'; I='17.74665496658379'; M='15.77454495387071'; N='16.11645684205958' }
    15 = @{ B='14.25824226632078';  C='11.28487288028479'; E='12.21014382835036'; F='21.29868154950795'; G='27.54500601051613'; H='12.90789535164634'; I='17.75168798116478'; M='15.73655255020567'; N='16.12388757623073' }
    16 = @{ B='13.87861842527938';  C='10.96217695221297'; E='12.11145130618689'; F='20.62722412089977'; G='27.24292495525691'; H='12.91181821877728'; I='17.78231003882133'; M='15.51810874379393'; N='16.1670308685499'  }
    17 = @{ B='13.64080227681377';  C='10.75907639443625'; E='12.0516124993018';  F='20.20408069597325'; G='27.05960827501367'; H='12.91521460255177'; I='17.8026772023528';  M='15.38352680775872'; N='16.19399940753076' }
    18 = @{ B='13.50223117666511';  C='10.64037205532478'; E='12.01746557029138'; F='19.95656407809801'; G='26.95495009209957'; H='12.91753106949776'; I='17.81497036666808'; M='15.30593751140416'; N='16.20969576544928' }
    19 = @{ B='13.45500934930159';  C='10.59985702749344'; E='12.0059515816088';  F='19.87204792380568'; G='26.91965303887292'; H='12.91837764311746'; I='17.81923172190996'; M='15.27963882127883'; N='16.21504206026318' }
    20 = @{ B='13.66630362339897';  C='10.78089205997381'; E='12.05795467367352'; F='20.24955283636154'; G='27.07904287422276'; H='12.9148154681148';  I='17.80044914797328'; M='15.39787272652501'; N='16.19110944810117' }
    21 = @{ B='14.35471012351132';  C='11.36659750344184'; E='12.23583318558064'; F='21.46857628470577'; G='27.62357857240974'; H='12.90718917584347'; I='17.74427307728919'; M='15.79274951744906'; N='16.11290350636678' }
    22 = @{ B='14.78820700252813';  C='11.73256700854454'; E='12.35425530504757'; F='22.22866616901552'; G='27.98544697802263'; H='12.90545315879856'; I='17.71276682960477'; M='16.04859469969713'; N='16.06342704191465' }
    23 = @{ B='14.55836596116395';  C='11.53878213169195'; E='12.29086465954454'; F='21.82633154458858'; G='27.79181042020741'; H='12.90608179033006'; I='17.72910309581263'; M='15.91227807501343'; N='16.08968436687815' }
    24 = @{ B='13.65478021283913';  C='10.77103523732007'; E='12.0550865790485';  F='20.22900810905287'; G='27.07025419772975'; H='12.91499478356262'; I='17.80145463335266'; M='15.39138760353448'; N='16.19241540119024' }
    25 = @{ B='12.6123103525052';   C='9.87082949385424';  E='11.81124908723674'; F='18.34778573295695'; G='26.32234962358157'; H='12.93847193250701'; I='17.90157863819771'; M='14.82300056306705'; N='16.31029770361473' }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = [double]$cols[$col]
    }
}
